$wb = $excel.ActiveWorkbook

# --- Sheet: y_fitted_on_begin_2016 ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Rows("2").Insert()
$ws1.Cells.Item(2,1).Value = 1991
$ws1.Cells.Item(2,2).Value = 68.02329931382778
$ws1.Cells.Item(3,1).Value = 1992
$ws1.Cells.Item(3,2).Value = 68.06687632970952
$ws1.Cells.Item(4,1).Value = 1993
$ws1.Cells.Item(4,2).Value = 68.5438299901289
$ws1.Cells.Item(5,1).Value = 1994
$ws1.Cells.Item(5,2).Value = 68.98772240466559
$ws1.Cells.Item(6,1).Value = 1995
$ws1.Cells.Item(6,2).Value = 68.91545451950354
$ws1.Cells.Item(7,1).Value = 1996
$ws1.Cells.Item(7,2).Value = 65.63193534547509
$ws1.Cells.Item(8,1).Value = 1997
$ws1.Cells.Item(8,2).Value = 65.78743840151026
$ws1.Cells.Item(9,1).Value = 1998
$ws1.Cells.Item(9,2).Value = 66.20082629547714
$ws1.Cells.Item(10,1).Value = 1999
$ws1.Cells.Item(10,2).Value = 66.57872160372045
$ws1.Cells.Item(11,1).Value = 2000
$ws1.Cells.Item(11,2).Value = 66.68904513991515
$ws1.Cells.Item(12,1).Value = 2001
$ws1.Cells.Item(12,2).Value = 66.6152315027792
$ws1.Cells.Item(13,1).Value = 2002
$ws1.Cells.Item(13,2).Value = 67.61442617075112
$ws1.Cells.Item(14,1).Value = 2003
$ws1.Cells.Item(14,2).Value = 68.43528711101776
$ws1.Cells.Item(15,1).Value = 2004
$ws1.Cells.Item(15,2).Value = 69.03396865700645
$ws1.Cells.Item(16,1).Value = 2005
$ws1.Cells.Item(16,2).Value = 69.43651585756811
$ws1.Cells.Item(17,1).Value = 2006
$ws1.Cells.Item(17,2).Value = 69.77072117389653
$ws1.Cells.Item(18,1).Value = 2007
$ws1.Cells.Item(18,2).Value = 69.65556947375694
$ws1.Cells.Item(19,1).Value = 2008
$ws1.Cells.Item(19,2).Value = 70.10021158000089
$ws1.Cells.Item(20,1).Value = 2009
$ws1.Cells.Item(20,2).Value = 70.25233692849474
$ws1.Cells.Item(21,1).Value = 2010
$ws1.Cells.Item(21,2).Value = 71.36444301958923
$ws1.Cells.Item(22,1).Value = 2011
$ws1.Cells.Item(22,2).Value = 70.76077494770277
$ws1.Cells.Item(23,1).Value = 2012
$ws1.Cells.Item(23,2).Value = 70.32317768917929
$ws1.Cells.Item(24,1).Value = 2013
$ws1.Cells.Item(24,2).Value = 70.30755671649591
$ws1.Cells.Item(25,1).Value = 2014
$ws1.Cells.Item(25,2).Value = 70.00729207014109
$ws1.Cells.Item(26,1).Value = 2015
$ws1.Cells.Item(26,2).Value = 70.04265577216293
$ws1.Cells.Item(27,1).Value = 2016
$ws1.Cells.Item(27,2).Value = 70.11009934902495

# --- Sheet: y_pred_on_2017_2021 ---
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws2.Cells.Item(2,1).Value = 2017
$ws2.Cells.Item(2,2).Value = 70.42247339660403
$ws2.Cells.Item(3,1).Value = 2018
$ws2.Cells.Item(3,2).Value = 69.89048660272515
$ws2.Cells.Item(4,1).Value = 2019
$ws2.Cells.Item(4,2).Value = 69.40702931894386
$ws2.Cells.Item(5,1).Value = 2020
$ws2.Cells.Item(5,2).Value = 68.94396714117777
$ws2.Cells.Item(6,1).Value = 2021
$ws2.Cells.Item(6,2).Value = 68.48294134922081

# --- Sheet: y_fitted_on_begin_2021 ---
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Rows("33").Delete()
$ws3.Cells.Item(2,1).Value = 1991
$ws3.Cells.Item(2,2).Value = 67.82358739635596
$ws3.Cells.Item(3,1).Value = 1992
$ws3.Cells.Item(3,2).Value = 67.94149549358059
$ws3.Cells.Item(4,1).Value = 1993
$ws3.Cells.Item(4,2).Value = 68.51930638285853
$ws3.Cells.Item(5,1).Value = 1994
$ws3.Cells.Item(5,2).Value = 69.06384947595568
$ws3.Cells.Item(6,1).Value = 1995
$ws3.Cells.Item(6,2).Value = 69.05750457861492
$ws3.Cells.Item(7,1).Value = 1996
$ws3.Cells.Item(7,2).Value = 65.65287872803188
$ws3.Cells.Item(8,1).Value = 1997
$ws3.Cells.Item(8,2).Value = 65.82200886714503
$ws3.Cells.Item(9,1).Value = 1998
$ws3.Cells.Item(9,2).Value = 66.22826625794823
$ws3.Cells.Item(10,1).Value = 1999
$ws3.Cells.Item(10,2).Value = 66.60685806181803
$ws3.Cells.Item(11,1).Value = 2000
$ws3.Cells.Item(11,2).Value = 66.75471779077789
$ws3.Cells.Item(12,1).Value = 2001
$ws3.Cells.Item(12,2).Value = 66.76016788316326
$ws3.Cells.Item(13,1).Value = 2002
$ws3.Cells.Item(13,2).Value = 67.72770514734015
$ws3.Cells.Item(14,1).Value = 2003
$ws3.Cells.Item(14,2).Value = 68.38590870003257
$ws3.Cells.Item(15,1).Value = 2004
$ws3.Cells.Item(15,2).Value = 68.92235253720612
$ws3.Cells.Item(16,1).Value = 2005
$ws3.Cells.Item(16,2).Value = 69.31000512660586
$ws3.Cells.Item(17,1).Value = 2006
$ws3.Cells.Item(17,2).Value = 69.61983347924664
$ws3.Cells.Item(18,1).Value = 2007
$ws3.Cells.Item(18,2).Value = 69.52447190167825
$ws3.Cells.Item(19,1).Value = 2008
$ws3.Cells.Item(19,2).Value = 69.97246213508186
$ws3.Cells.Item(20,1).Value = 2009
$ws3.Cells.Item(20,2).Value = 70.12337658059307
$ws3.Cells.Item(21,1).Value = 2010
$ws3.Cells.Item(21,2).Value = 71.30111256491496
$ws3.Cells.Item(22,1).Value = 2011
$ws3.Cells.Item(22,2).Value = 70.71595366274323
$ws3.Cells.Item(23,1).Value = 2012
$ws3.Cells.Item(23,2).Value = 70.39652766148542
$ws3.Cells.Item(24,1).Value = 2013
$ws3.Cells.Item(24,2).Value = 70.5200770629229
$ws3.Cells.Item(25,1).Value = 2014
$ws3.Cells.Item(25,2).Value = 70.33624594338299
$ws3.Cells.Item(26,1).Value = 2015
$ws3.Cells.Item(26,2).Value = 70.49783487552872
$ws3.Cells.Item(27,1).Value = 2016
$ws3.Cells.Item(27,2).Value = 70.70148148665541
$ws3.Cells.Item(28,1).Value = 2017
$ws3.Cells.Item(28,2).Value = 71.10953125770007
$ws3.Cells.Item(29,1).Value = 2018
$ws3.Cells.Item(29,2).Value = 70.88552500370633
$ws3.Cells.Item(30,1).Value = 2019
$ws3.Cells.Item(30,2).Value = 70.8837952589852
$ws3.Cells.Item(31,1).Value = 2020
$ws3.Cells.Item(31,2).Value = 70.91173969906353
$ws3.Cells.Item(32,1).Value = 2021
$ws3.Cells.Item(32,2).Value = 72.10773953836764

# --- Sheet: y_pred_on_2022_2026 ---
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws4.Cells.Item(2,1).Value = 2022
$ws4.Cells.Item(2,2).Value = 71.45030869860169
$ws4.Cells.Item(3,1).Value = 2023
$ws4.Cells.Item(3,2).Value = 71.30861486680428
$ws4.Cells.Item(4,1).Value = 2024
$ws4.Cells.Item(4,2).Value = 71.19701074669356
$ws4.Cells.Item(5,1).Value = 2025
$ws4.Cells.Item(5,2).Value = 71.10738866087041
$ws4.Cells.Item(6,1).Value = 2026
$ws4.Cells.Item(6,2).Value = 71.03295749773075
